$wb = $excel.ActiveWorkbook

# Rename "f__Lachnospiracea-t" -> "f__Lachnospiraceae-t"
$renameSheet = $wb.Worksheets.Item("f__Lachnospiracea-t")
$renameSheet.Name = "f__Lachnospiraceae-t"

# Delete "g__CAG-791-t"
$excel.DisplayAlerts = $false
$deleteSheet = $wb.Worksheets.Item("g__CAG-791-t")
$deleteSheet.Delete()
$excel.DisplayAlerts = $true
